$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# Update the feature-count value in B9 (KNN numeric results row)
$ws.Range("B9").Value = 1766

# Add new KNN numeric results in H9:M9 as text values (reusing 0.258 already in shared strings).
# Build each value with a TEXT() formula (so Excel does not auto-convert the
# numeric-looking string into a real number), then flatten the formula down to
# a plain value in-place. This keeps the cell's existing style (inherited from
# the column, xf index 3) untouched and avoids creating any new number-format
# style entries.
$vals = @("0.512", "0.270", "0.512", "0.269", "0.512", "0.258")
$cols = @("H", "I", "J", "K", "L", "M")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "9")
    $cell.Formula = '=TEXT(' + $vals[$i] + ',"0.000")'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}
$excel.CutCopyMode = $false

# Leave the selection where the author ended up after entering the new results
$ws.Range("M10").Select() | Out-Null
